# "discussion - minor issues fixed"
#
# Functional changes captured from the diff:
#   1. C2:C5  10 -> 5      (run_num column)
#   2. D2:D5  20 -> 10     (period_num column)
#   3. I1:I5  drop the stray "applyAlignment" cell style (s=4 -> s=2),
#      i.e. the influence_msg column gets the same plain style used by
#      the rest of the header/body text cells.
#   4. Selection moves to D2 (and the view scrolls so column E is the
#      left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: update the run_num / period_num values for rows 2-5 -----------
$ws.Range("C2:C5").Value = 5
$ws.Range("D2:D5").Value = 10

# --- 3: re-style I1:I5 to match the other "s=2" cells (e.g. E1) -----------
# Copy/PasteSpecial(formats) reuses the existing cellXf for that font
# instead of minting a new one, so I1:I5 land back on the same style index
# used elsewhere in the sheet.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("I1:I5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 4: selection / scroll position ---------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("D2").Select() | Out-Null
